$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update region names (column A) for rows whose region changed
$ws.Range("A2").Value = "Roraima"
$ws.Range("A3").Value = "Espírito Santo"
$ws.Range("A4").Value = "Maranhão"
$ws.Range("A5").Value = "Rio Grande do Norte"
$ws.Range("A6").Value = "Bahia"
$ws.Range("A7").Value = "Tocantins"

# Update period label (column B) for rows 2-10 from 2025/04-2024/04 to 2025/07-2024/07
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = "Diferença 2025/07 - 2024/07"
}

# Update values (column C)
$ws.Range("C2").Value = 1.62
$ws.Range("C3").Value = 1.55
$ws.Range("C4").Value = 1.53
$ws.Range("C5").Value = 1.39
$ws.Range("C6").Value = 1.2
$ws.Range("C7").Value = 1.16
$ws.Range("C8").Value = 0.7
$ws.Range("C9").Value = 0.78
$ws.Range("C10").Value = 0.85

# Update ranking (column D) for row 8 only
$ws.Range("D8").Value = "14º"
